$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seguimiento grado 10 - actualizacion al 27/04/2015
#
# The order in which shared-string text is assigned matters: this engine
# rebuilds the shared-string table from scratch on save, dropping entries
# that lose their last referrer and appending freshly-seen text in the
# order it is encountered. Driving the edits in the same order the final
# workbook's string table is laid out reproduces the exact table used by
# the real edit.

# Row 6 / Row 7: "H" column note is rewritten (same new text in both
# cells so the shared string is replaced in place instead of forked).
$ws.Cells.Item(6, 8).Value = "Se van a incluir nuevas imágenes"
$ws.Cells.Item(7, 8).Value = "Se van a incluir nuevas imágenes"

# Row 7: a review date is filled in. Set the value first, then copy the
# date format from G6 (style already has the date number format) so the
# cell picks up the existing date style instead of minting a new one.
$ws.Cells.Item(7, 7).Value = 42078
$ws.Cells.Item(6, 7).Copy()
$ws.Cells.Item(7, 7).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 7 & 13 share the same "F" note - update both to the same text so
# it is replaced in place rather than split into two entries.
$ws.Cells.Item(7, 6).Value = "Está siendo reestructurada"
$ws.Cells.Item(13, 6).Value = "Está siendo reestructurada"

# Rows 14 then 8 introduce two brand-new notes (in this order, matching
# the order the original edit appended them to the shared-string table).
$ws.Cells.Item(14, 6).Value = "Entregada por la autora para edición"
$ws.Cells.Item(8, 6).Value = "Entregada por el autor para edición"

# Row 8: manuscript/delivery dates.
$ws.Cells.Item(8, 2).Value = 42118
$ws.Cells.Item(8, 3).Value = 42118

# Row 14: manuscript/delivery dates.
$ws.Cells.Item(14, 2).Value = 42121
$ws.Cells.Item(14, 3).Value = 42121

# Row heights follow the new wrapped-text content.
$ws.Rows.Item(6).RowHeight = 29.25
$ws.Rows.Item(8).RowHeight = 29.25
$ws.Rows.Item(14).RowHeight = 29.25

# Leave the selection where the author finished working.
[void]$ws.Range("F15").Select()
